$p = $ppt.ActivePresentation
Write-Host "no-op"
